$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new task row (row 22) at the bottom of the list, reusing the
# formatting of an existing "Cделано"/"0.6.0" row (row 7) so the new
# row's cell styles match what a user would get by duplicating a row.
$ws.Range("A7:C7").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A22").Value = "Сделать возможность бить руками"
$ws.Range("B22").Value = "Cделано"
$ws.Range("C22").Value = "0.6.0"

# Update the current selection to reflect where the user ended up editing.
$ws.Range("D14").Select()
